$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Columns("I:I").Insert()
$ws.Range("I2").Value = "GameSettings"
$ws.Range("I1").Clear()
$ws.Columns("I:I").ColumnWidth = 18.666666666666668
